# Applies the "minor language changes" edit described in the commit.
# Targets the three USB-whitelisting slides (Gov_templates/USB_whitelist.pptx).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: "SMART INDIA HACKATHON '18"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# --- Content Placeholder 2 (shape 2): tidy up tab runs -----------------
$cp = $s1.Shapes.Item(2)
$cpTr = $cp.TextFrame.TextRange

# Paragraph 2 ("Problem Statement: ...intranet") - merge trailing tabs into
# the preceding run so the paragraph ends as a single run.
$para2 = $cpTr.Paragraphs(2)
$r2 = $para2.Runs(2)
$r2.Text = ": Prototype/application for whitelisting of USB devices in OFB which can be subsequently used on internet as well as on intranet" + "				"
$para2.Runs(3).Text = ""

# Paragraph 3 ("Problem Code: ... Team Name:") - merge "Problem " + "Code"
$para3 = $cpTr.Paragraphs(3)
$para3.Runs(1).Text = "Problem Code"
$para3.Runs(2).Text = ""

# Paragraph 4 ("Team Leader Name: ... College Code:") - merge the two
# single-tab runs into a single double-tab run.
$para4 = $cpTr.Paragraphs(4)
$para4.Runs(6).Text = "		"
$para4.Runs(7).Text = ""

# --- Rectangle 5 (shape 4): grow the box + rewrite the bullet copy ------
$rect = $s1.Shapes.Item(4)
$rect.Height = 4247317 / 12700
$rectTr = $rect.TextFrame.TextRange

$rectTr.Paragraphs(2).Runs(1).Text = "This can be implemented using Blockchain system which provides tamper-proof dataset e.g. – Hyperledger. We will implement encryption on top of blockchain infrastructure"

$rectTr.Paragraphs(3).Runs(1).Text = "When the storage device is connected to the computer a program will extract the MAC address and check  if it is present on the database. A local probabilistic database (e.g. like bloom filter) will be used to check if the MAC address is while listed. If  MAC address is not on the white list then the OS event will notify the user regarding the same and block it. Use of Probabilistic database will allow the system to work on internet/intranet or disconnected computers from the unauthorized storage devices. The program used, handles the connectivity of storage devices with the computer."
$rectTr.Paragraphs(3).Runs(2).Text = ""

$rectTr.Paragraphs(4).Runs(1).Text = "Another level of security will be creating  an encrypted filesystem on storage device. Filesystem will get decrypted automatically when a whitelisted storage device is connected to an authorized computer. If decrypted correctly the data transfer or access between storage device and computer is possible."

# Drop the trailing empty bullet paragraph.
$rectTr.Paragraphs(5).Delete()

# ---------------------------------------------------------------------------
# Slide 2: "TECHNOLOGY STACK" / "DEPENDENCIES / SHOW STOPPER"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# --- TextBox 5 (shape 1): how the solution works ------------------------
$tb5 = $s2.Shapes.Item(1)
$tb5Tr = $tb5.TextFrame.TextRange

$tb5Tr.Paragraphs(1).Runs(1).Text = "The encrypted filesystem will be created on the storage device when it is connected to the blockchain system for the first time for registration on the database. "
$tb5Tr.Paragraphs(2).Runs(1).Text = "The computers which are offline  will use a probabilistic database of authorized MAC addresses (e.g. like Bloom filter). The driver program on the computer will check the storage device against this probabilistic database. The driver program will block the devices which are not on the whitelist. "

# Paragraph 5 becomes the new (and last) bullet; delete paragraphs 3 and 4.
$tb5Tr.Paragraphs(5).Runs(1).Text = "The Bloom filter database will be updated when it is connected to the internet periodically."
$tb5Tr.Paragraphs(4).Delete()
$tb5Tr.Paragraphs(3).Delete()

# --- TextBox 3 (shape 3): technology stack bullets ----------------------
$tb3 = $s2.Shapes.Item(3)
$tb3Tr = $tb3.TextFrame.TextRange

$tb3Tr.Paragraphs(1).Runs(1).Text = "Opensource Blockchain distributed database  like Hyperle"
$tb3Tr.Paragraphs(2).Runs(1).Text = "Web application based on blockchain database"
$tb3Tr.Paragraphs(4).Runs(1).Text = "Windows service for detecting/blocking  connected USB storage device"

# --- Rectangle 6 (shape 5): dependencies / show stoppers ----------------
$rect6 = $s2.Shapes.Item(5)
$rect6Tr = $rect6.TextFrame.TextRange

$rect6Tr.Paragraphs(1).Runs(1).Text = "Probabilistic data structures are not 100% accurate. In extremely rare case authorized device may get blocked if computer is offline"
$rect6Tr.Paragraphs(2).Runs(1).Text = "The offline devices have to be connected to the internet or intranet for updating it periodicatlly. This can give access to the devices which have been removed from the whitelist and won't give access to the newly added devices."

# Drop the old third bullet entirely.
$rect6Tr.Paragraphs(3).Delete()

Write-Output "Edit complete"
